$wb = $excel.ActiveWorkbook

# --- Rename sheets (date range moved from July to September) ---
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws1.Name = "9-2-13"
$ws2.Name = "9-9-13"

# --- Sheet2 ("9-9-13") data updates ---

# Row 1: due date moves to new format/style, C1 becomes empty
$ws2.Range("B1").NumberFormat = "mm/dd/yy"
$ws2.Range("B1").Value = 40059
$ws2.Range("C1").Clear()

# Row 3 header cells get the new date style (text stays the same)
$ws2.Range("B3").NumberFormat = "mm/dd/yy"
$ws2.Range("C3").NumberFormat = "mm/dd/yy"

# Row 4
$ws2.Range("B4").NumberFormat = "mm/dd/yy"
$ws2.Range("B4").Value = 40059
$ws2.Range("C4").Clear()
$ws2.Range("M4").Value = 39999
$ws2.Range("N4").Value = 0.5

# Row 5: B5 currently holds the shared string "?" (index 27, text changed below);
# repoint it to a date value and restyle, same for C5
$ws2.Range("B5").Value = "Set up Vex Dev Environment"
$ws2.Range("A6").Value = "Set up Vex Dev Environment"
$ws2.Range("B5").NumberFormat = "mm/dd/yy"
$ws2.Range("B5").Value = 40061
$ws2.Range("C5").NumberFormat = "mm/dd/yy"
$ws2.Range("C5").Value = 40058

# Row 6 (new)
$ws2.Range("B6").NumberFormat = "mm/dd/yy"
$ws2.Range("B6").Value = 40059
$ws2.Range("C6").NumberFormat = "mm/dd/yy"
$ws2.Range("C6").Value = 40062
$ws2.Range("D6").NumberFormat = "0.0%"
$ws2.Range("D6").Value = 1
$ws2.Range("E6").Value = 9

# Row 7 (new)
$ws2.Range("A7").Value = "Find motor spec Sheets"
$ws2.Range("B7").NumberFormat = "mm/dd/yy"
$ws2.Range("B7").Value = 40063
$ws2.Range("C7").NumberFormat = "mm/dd/yy"
$ws2.Range("C7").Value = 40062
$ws2.Range("D7").NumberFormat = "0.0%"
$ws2.Range("D7").Value = 1
$ws2.Range("E7").Value = 0.5
$ws2.Range("F7").Value = "Start 2:30"

# Row 8 (new)
$ws2.Range("A8").Value = "Code individual Motor motion"
$ws2.Range("B8").NumberFormat = "mm/dd/yy"
$ws2.Range("B8").Value = 40063
$ws2.Range("C8").NumberFormat = "mm/dd/yy"
$ws2.Range("C8").Value = 40062
$ws2.Range("D8").NumberFormat = "0.0%"
$ws2.Range("D8").Value = 1
$ws2.Range("E8").Value = 0.25
$ws2.Range("F8").Value = "Start 3:00"

# Row 9 (new)
$ws2.Range("A9").Value = "Test individual Motor Motion"
$ws2.Range("B9").NumberFormat = "mm/dd/yy"
$ws2.Range("B9").Value = 40063
$ws2.Range("C9").NumberFormat = "mm/dd/yy"
$ws2.Range("C9").Value = 40062
$ws2.Range("D9").NumberFormat = "0.0%"
$ws2.Range("D9").Value = 1
$ws2.Range("E9").Value = 0.25

# Row 10 (new)
$ws2.Range("A10").Value = "Replace 3-wire motors with 2-wire"
$ws2.Range("B10").NumberFormat = "mm/dd/yy"
$ws2.Range("B10").Value = 40061
$ws2.Range("C10").NumberFormat = "mm/dd/yy"
$ws2.Range("C10").Value = 40062
$ws2.Range("D10").NumberFormat = "0.0%"
$ws2.Range("D10").Value = 1
$ws2.Range("E10").Value = 0.75
$ws2.Range("F10").Value = "Start 3:30"

# Row 11 (new)
$ws2.Range("A11").Value = "Code Forward/Backward/Left/Right"
$ws2.Range("B11").NumberFormat = "mm/dd/yy"
$ws2.Range("B11").Value = 40063
$ws2.Range("C11").NumberFormat = "mm/dd/yy"
$ws2.Range("C11").Value = 40062
$ws2.Range("D11").NumberFormat = "0.0%"
$ws2.Range("D11").Value = 1
$ws2.Range("E11").Value = 0.75

# Row 12 (new)
$ws2.Range("A12").Value = "Test Forward/Backward/Left/Right"
$ws2.Range("B12").NumberFormat = "mm/dd/yy"
$ws2.Range("B12").Value = 40063
$ws2.Range("C12").NumberFormat = "mm/dd/yy"
$ws2.Range("C12").Value = 40062
$ws2.Range("D12").NumberFormat = "0.0%"
$ws2.Range("D12").Value = 1
$ws2.Range("E12").Value = 0.25

# Row 13 (new)
$ws2.Range("A13").Value = "Read documentation on Vex Controller Communication"
$ws2.Range("B13").NumberFormat = "mm/dd/yy"
$ws2.Range("B13").Value = 40063
$ws2.Range("D13").NumberFormat = "0.0%"
$ws2.Range("D13").Value = 0.05
$ws2.Range("E13").Value = 0.5

# --- Column B:C formatting/width (cosmetic, matches the new date column style) ---
$ws2.Columns("B:C").ColumnWidth = 10

# --- Selection / view state ---
$ws2.Range("J8").Select()
